$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.823.39'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '2.929.95'
$ws.Range('E3').Value = '  +3.78%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '352.78'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '111.83'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.69%  '
$ws.Range('E7').Value = '  +0.50%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.622'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.35'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.75%  '
$ws.Range('E11').Value = '  +5.31%  '
$ws.Range('E12').Value = '  +1.21%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.97'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.30%  '
$ws.Range('D14').Value = '3.395.04'
$ws.Range('E14').Value = '  +3.80%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.73'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.19%  '
$ws.Range('D16').Value = '2.934.35'
$ws.Range('E16').Value = '  +4.18%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.984'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.52%  '
$ws.Range('D18').Value = '51.896.36'
$ws.Range('E18').Value = '  +0.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.31'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.62'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.17'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +6.00%  '
$ws.Range('E22').Value = '  +1.69%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '71.24'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '268.43'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.30%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.78'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.94%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.180'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +10.83%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '26.89'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.75%  '
$ws.Range('E28').Value = '  +0.34%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.29'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +16.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.105'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +17.67%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '10.55'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '37.26'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.70%  '
$ws.Range('B34').Value = 'RenderToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.22'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +10.67%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '52.82'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('E36').Value = '  +1.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.18%  '
$ws.Range('E38').Value = '  +3.89%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.86'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.42%  '
$ws.Range('E40').Value = '  +1.90%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.69'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.11%  '
$ws.Range('E42').Value = '  +1.40%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '23.17'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.06%  '
$ws.Range('E44').Value = '  -1.04%  '
$ws.Range('E45').Value = '  +1.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.51'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.63%  '
$ws.Range('D47').Value = '2.168.65'
$ws.Range('E47').Value = '  -0.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '111.59'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -9.30%  '
$ws.Range('E49').Value = '  -0.18%  '
$ws.Range('E50').Value = '  +12.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.939'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.61%  '
